# Tottus: bump existing rows 2-3 to the new order/date, and append the
# remaining lines (rows 4-16) of the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Common (shared) values for every row of this order.
$nroOrden   = "47035710"
$rut        = "78627210-6"
$nombre     = "Hipermercados TOTTUS SA"
$sucursal   = "CD TOTTUS"
$fechaEmision    = "20250814"
$fechaCompromiso = "20250826"
$nroOrdenSalida  = "7637"

# Per-line SKU (col S) + quantity (col T), row 2 through row 16.
$lines = @{
    2  = @{ S = "20215631"; T = 24  }
    3  = @{ S = "20287252"; T = 24  }
    4  = @{ S = "20287256"; T = 48  }
    5  = @{ S = "20287253"; T = 24  }
    6  = @{ S = "20202318"; T = 24  }
    7  = @{ S = "20215632"; T = 240 }
    8  = @{ S = "20202336"; T = 24  }
    9  = @{ S = "20202335"; T = 12  }
    10 = @{ S = "20202339"; T = 72  }
    11 = @{ S = "20202309"; T = 60  }
    12 = @{ S = "20202351"; T = 24  }
    13 = @{ S = "20202337"; T = 12  }
    14 = @{ S = "20287251"; T = 24  }
    15 = @{ S = "20202310"; T = 12  }
    16 = @{ S = "20215634"; T = 48  }
}

foreach ($r in ($lines.Keys | Sort-Object)) {
    $row = $lines[$r]

    # Text-valued columns must stay text (not get re-interpreted as numbers).
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $nroOrden
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $nroOrden

    if ($r -ge 4) {
        $ws.Range("C$r").NumberFormat = "@"
        $ws.Range("C$r").Value = $rut
        $ws.Range("D$r").Value = $nombre
        $ws.Range("E$r").Value = $sucursal
        $ws.Range("F$r").Value = $sucursal
    }

    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = $fechaEmision
    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value = $fechaCompromiso

    $ws.Range("S$r").NumberFormat = "@"
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T

    $ws.Range("AB$r").NumberFormat = "@"
    $ws.Range("AB$r").Value = $nroOrdenSalida
}

Write-Output "tottus rows updated/added"
